# Confusion Matrix Excel File
# - Rename the "Purchase" label to "Will Purchase" everywhere it is used
#   (B4, B5, B8, B9, C3, C5, C7, C9) so the shared string itself is updated
#   in place rather than creating a brand-new string entry.
# - Move the active selection from F13 to G15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Will Purchase"
$ws.Range("B4").Value = "Will Purchase"
$ws.Range("B5").Value = "Will Purchase"
$ws.Range("C5").Value = "Will Purchase"
$ws.Range("C7").Value = "Will Purchase"
$ws.Range("B8").Value = "Will Purchase"
$ws.Range("B9").Value = "Will Purchase"
$ws.Range("C9").Value = "Will Purchase"

$ws.Range("G15").Select()
